# Applies the "Fruta / hortaliza, semanal" update to the Camote sheet:
# shifts existing weekly records and appends newly scraped rows, extending
# the data range from A1:R33 to A1:R37.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("D8").Value = 44424
$ws.Range("J8").Value = 790
$ws.Range("K8").Value = 13000
$ws.Range("L8").Value = 14000
$ws.Range("M8").Value = 13506
$ws.Range("N8").Value = '$/malla 18 kilos'
$ws.Range("P8").Value = 750
$ws.Range("Q8").Value = 18

# Row 9
$ws.Range("D9").Value = 44424
$ws.Range("I9").Value = 'Segunda'
$ws.Range("J9").Value = 520
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 12000
$ws.Range("P9").Value = 667

# Row 10
$ws.Range("D10").Value = 44403
$ws.Range("J10").Value = 1330
$ws.Range("K10").Value = 11000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 11500
$ws.Range("N10").Value = '$/caja 15 kilos granel'
$ws.Range("P10").Value = 767
$ws.Range("Q10").Value = 15

# Row 11
$ws.Range("D11").Value = 44396
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 1330
$ws.Range("K11").Value = 9000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 9500
$ws.Range("P11").Value = 528

# Row 12
$ws.Range("D12").Value = 44340
$ws.Range("J12").Value = 1420
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 11000
$ws.Range("M12").Value = 10500
$ws.Range("P12").Value = 583

# Row 13
$ws.Range("D13").Value = 44340
$ws.Range("I13").Value = 'Segunda'
$ws.Range("J13").Value = 970
$ws.Range("K13").Value = 8000
$ws.Range("L13").Value = 8000
$ws.Range("M13").Value = 8000
$ws.Range("P13").Value = 444

# Row 14
$ws.Range("D14").Value = 44221
$ws.Range("L14").Value = 12000
$ws.Range("M14").Value = 11444
$ws.Range("P14").Value = 636

# Row 15
$ws.Range("D15").Value = 44242
$ws.Range("K15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("M15").Value = 10000
$ws.Range("P15").Value = 556

# Row 16
$ws.Range("D16").Value = 44200
$ws.Range("J16").Value = 1800
$ws.Range("K16").Value = 11000
$ws.Range("L16").Value = 11000
$ws.Range("M16").Value = 11000
$ws.Range("P16").Value = 611

# Row 17
$ws.Range("D17").Value = 44298
$ws.Range("J17").Value = 1600
$ws.Range("L17").Value = 14000
$ws.Range("M17").Value = 14000
$ws.Range("P17").Value = 778

# Row 18
$ws.Range("D18").Value = 44354
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 700
$ws.Range("K18").Value = 14000
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = 14500
$ws.Range("P18").Value = 806

# Row 19
$ws.Range("D19").Value = 44410
$ws.Range("J19").Value = 970
$ws.Range("K19").Value = 14000
$ws.Range("L19").Value = 15000
$ws.Range("M19").Value = 14505
$ws.Range("P19").Value = 806

# Row 20
$ws.Range("D20").Value = 44410
$ws.Range("I20").Value = 'Segunda'
$ws.Range("J20").Value = 340
$ws.Range("L20").Value = 12000
$ws.Range("M20").Value = 12000
$ws.Range("P20").Value = 667

# Row 21
$ws.Range("D21").Value = 44172
$ws.Range("J21").Value = 1600
$ws.Range("K21").Value = 12000
$ws.Range("L21").Value = 12000
$ws.Range("M21").Value = 12000
$ws.Range("P21").Value = 667

# Row 22
$ws.Range("D22").Value = 44214
$ws.Range("J22").Value = 1900
$ws.Range("K22").Value = 12000
$ws.Range("L22").Value = 13000
$ws.Range("M22").Value = 12526
$ws.Range("P22").Value = 696

# Row 23
$ws.Range("D23").Value = 44193
$ws.Range("J23").Value = 1800
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 11000
$ws.Range("M23").Value = 10556
$ws.Range("P23").Value = 586

# Row 24
$ws.Range("D24").Value = 44389
$ws.Range("J24").Value = 1420
$ws.Range("K24").Value = 8000
$ws.Range("L24").Value = 9000
$ws.Range("M24").Value = 8500
$ws.Range("P24").Value = 472

# Row 25
$ws.Range("D25").Value = 44305
$ws.Range("K25").Value = 12000
$ws.Range("L25").Value = 12000
$ws.Range("M25").Value = 12000
$ws.Range("P25").Value = 667

# Row 26
$ws.Range("D26").Value = 44417
$ws.Range("J26").Value = 880
$ws.Range("K26").Value = 14000
$ws.Range("L26").Value = 15000
$ws.Range("M26").Value = 14500
$ws.Range("P26").Value = 806

# Row 27
$ws.Range("D27").Value = 44417
$ws.Range("I27").Value = 'Segunda'
$ws.Range("J27").Value = 340
$ws.Range("K27").Value = 12000
$ws.Range("M27").Value = 12000
$ws.Range("P27").Value = 667

# Row 28
$ws.Range("D28").Value = 44333
$ws.Range("J28").Value = 1410
$ws.Range("L28").Value = 11000
$ws.Range("M28").Value = 10500
$ws.Range("P28").Value = 583

# Row 29
$ws.Range("D29").Value = 44277
$ws.Range("J29").Value = 1600
$ws.Range("K29").Value = 15000
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = 15000
$ws.Range("P29").Value = 833

# Row 30
$ws.Range("D30").Value = 44319
$ws.Range("J30").Value = 1510

# Row 31
$ws.Range("D31").Value = 44186
$ws.Range("J31").Value = 1800
$ws.Range("K31").Value = 11000
$ws.Range("L31").Value = 12000
$ws.Range("M31").Value = 11556
$ws.Range("P31").Value = 642

# Row 32
$ws.Range("D32").Value = 44326
$ws.Range("J32").Value = 1600
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = 10000
$ws.Range("P32").Value = 556

# Row 33
$ws.Range("D33").Value = 44382
$ws.Range("J33").Value = 1510
$ws.Range("K33").Value = 8000
$ws.Range("L33").Value = 9000
$ws.Range("M33").Value = 8500
$ws.Range("P33").Value = 472

# Row 34
$ws.Range("A34").Value = 9
$ws.Range("B34").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C34").Value = 'Metropolitana'
$ws.Range("D34").Value = 44270
$ws.Range("D34").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E34").Value = 13
$ws.Range("F34").Value = 100114002
$ws.Range("G34").Value = 'Camote'
$ws.Range("H34").Value = 'Sin especificar'
$ws.Range("I34").Value = 'Primera'
$ws.Range("J34").Value = 16000
$ws.Range("K34").Value = 10000
$ws.Range("L34").Value = 11000
$ws.Range("M34").Value = 10500
$ws.Range("N34").Value = '$/malla 18 kilos'
$ws.Range("O34").Value = 'Perú'
$ws.Range("P34").Value = 583
$ws.Range("Q34").Value = 18
$ws.Range("R34").Value = 'Hortaliza'

# Row 35
$ws.Range("A35").Value = 9
$ws.Range("B35").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C35").Value = 'Metropolitana'
$ws.Range("D35").Value = 44376
$ws.Range("D35").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E35").Value = 13
$ws.Range("F35").Value = 100114002
$ws.Range("G35").Value = 'Camote'
$ws.Range("H35").Value = 'Sin especificar'
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 520
$ws.Range("K35").Value = 9000
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = 9500
$ws.Range("N35").Value = '$/malla 18 kilos'
$ws.Range("O35").Value = 'Perú'
$ws.Range("P35").Value = 528
$ws.Range("Q35").Value = 18
$ws.Range("R35").Value = 'Hortaliza'

# Row 36
$ws.Range("A36").Value = 9
$ws.Range("B36").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C36").Value = 'Metropolitana'
$ws.Range("D36").Value = 44179
$ws.Range("D36").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E36").Value = 13
$ws.Range("F36").Value = 100114002
$ws.Range("G36").Value = 'Camote'
$ws.Range("H36").Value = 'Sin especificar'
$ws.Range("I36").Value = 'Primera'
$ws.Range("J36").Value = 1500
$ws.Range("K36").Value = 10000
$ws.Range("L36").Value = 11000
$ws.Range("M36").Value = 10600
$ws.Range("N36").Value = '$/malla 18 kilos'
$ws.Range("O36").Value = 'Perú'
$ws.Range("P36").Value = 589
$ws.Range("Q36").Value = 18
$ws.Range("R36").Value = 'Hortaliza'

# Row 37
$ws.Range("A37").Value = 9
$ws.Range("B37").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C37").Value = 'Metropolitana'
$ws.Range("D37").Value = 44284
$ws.Range("D37").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E37").Value = 13
$ws.Range("F37").Value = 100114002
$ws.Range("G37").Value = 'Camote'
$ws.Range("H37").Value = 'Sin especificar'
$ws.Range("I37").Value = 'Primera'
$ws.Range("J37").Value = 1600
$ws.Range("K37").Value = 12000
$ws.Range("L37").Value = 12000
$ws.Range("M37").Value = 12000
$ws.Range("N37").Value = '$/malla 18 kilos'
$ws.Range("O37").Value = 'Perú'
$ws.Range("P37").Value = 667
$ws.Range("Q37").Value = 18
$ws.Range("R37").Value = 'Hortaliza'
